$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New account-statement data: previous periods removed, new periods (2107..2101 roughly)
# added for both workers. Column layout: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo,
# F=Valor Mora, G=Salario Basico

$data = @(
    @("CC", "72135702",   "JORGE MARIO LONDOÑO CARDONA",  "2107", 29260, 908526),
    @("CC", "72135702",   "JORGE MARIO LONDOÑO CARDONA",  "2106", 35112, 908526),
    @("CC", "72135702",   "JORGE MARIO LONDOÑO CARDONA",  "2105", 35112, 908526),
    @("CC", "72135702",   "JORGE MARIO LONDOÑO CARDONA",  "2104", 35112, 828116),
    @("CC", "72135702",   "JORGE MARIO LONDOÑO CARDONA",  "2103", 35112, 828116),
    @("CC", "72135702",   "JORGE MARIO LONDOÑO CARDONA",  "2102", 35112, 908526),
    @("CC", "72135702",   "JORGE MARIO LONDOÑO CARDONA",  "2003", 35112, 908526),
    @("CC", "72135702",   "JORGE MARIO LONDOÑO CARDONA",  "2001", 33125, 908526),
    @("CC", "1047484176", "YESICA PATRICIA LONDOÑO MEZA", "2107", 29260, 908526),
    @("CC", "1047484176", "YESICA PATRICIA LONDOÑO MEZA", "2106", 35112, 908526),
    @("CC", "1047484176", "YESICA PATRICIA LONDOÑO MEZA", "2105", 35112, 908526),
    @("CC", "1047484176", "YESICA PATRICIA LONDOÑO MEZA", "2104", 35112, 908526),
    @("CC", "1047484176", "YESICA PATRICIA LONDOÑO MEZA", "2103", 35112, 908526),
    @("CC", "1047484176", "YESICA PATRICIA LONDOÑO MEZA", "2102", 35112, 908526)
)

$row = 16
foreach ($rec in $data) {
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
    $row++
}
